# Rewrite Sheet1 in the new column layout:
#  - drop distance_km / match_score / match_rank
#  - keep the metadata columns (cand_gender .. job_work_province)
#  - reorder the per-language count columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("cand_gender", "cand_age_bucket", "cand_domicile_province", "cand_domicile_region", "job_contract_type", "job_work_province", "Svedese", "Spagnolo", "Finlandese", "Ebraico", "Persiano", "Portoghese", "Cinese", "Arabo", "Serbo", "Albanese", "Croato", "Ceco", "Danese", "Rumeno", "Macedone", "Tedesco")
    ,@(0, 0, 184, 11, 0, 185, 0, 0, 0, 0, 5, 0, 0, 6, 0, 0, 1, 0, 0, 0, 1, 0)
    ,@(0, 0, 106, 22, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0)
    ,@(0, 0, 585, 128, 0, 50, 0, 0, 0, 0, 0, 0, 1, 2, 0, 1, 0, 0, 0, 0, 0, 0)
    ,@(0, 0, 223, 180, 0, 85, 3, 0, 0, 0, 2, 0, 2, 1, 0, 0, 0, 0, 1, 0, 0, 0)
    ,@(0, 0, 953, 0, 0, 0, 2, 0, 0, 0, 1, 0, 0, 0, 0, 2, 0, 0, 0, 0, 1, 0)
    ,@(0, 0, 464, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 5, 0, 0, 0, 0, 0, 0)
    ,@(0, 0, 484, 159, 0, 180, 0, 0, 0, 0, 1, 0, 0, 1, 0, 2, 0, 0, 0, 0, 0, 0)
    ,@(0, 0, 1721, 1260, 0, 412, 0, 0, 2, 0, 0, 0, 0, 10, 0, 0, 0, 1, 0, 0, 2, 0)
    ,@(0, 0, 955, 248, 0, 780, 2, 0, 1, 0, 6, 1, 0, 0, 0, 4, 0, 0, 0, 0, 0, 0)
    ,@(0, 0, 346, 85, 0, 0, 0, 0, 0, 0, 1, 0, 0, 1, 0, 0, 4, 0, 0, 0, 0, 3)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Drop the 3 now-unused trailing columns (old width 25 -> new width 22)
$ws.Range("W1:Y11").Delete()
